$d = $word.ActiveDocument

# The document uses placeholder tokens [[PERSON_N]]. PERSON_9 was identified as
# a duplicate of PERSON_8 and removed; every PERSON_N with N >= 10 is shifted
# down by one (PERSON_10 -> PERSON_9, PERSON_11 -> PERSON_10, ...), and the
# former "od [[PERSON_9]]" slot becomes "od [[PERSON_8]]". Applying the
# replacements in ascending order of N means each Find.Execute pass only ever
# matches the still-untouched, higher-numbered original tokens, so the shifts
# do not cascade into each other.
for ($n = 9; $n -le 47; $n++) {
    $old = "[[PERSON_$n]]"
    $new = "[[PERSON_" + ($n - 1) + "]]"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
